$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 192
$ws1.Range("F3").Value = 109
$ws1.Range("F6").Value = 5298
$ws1.Range("F7").Value = 446
$ws1.Range("F8").Value = 633
$ws1.Range("F9").Value = 914
$ws1.Range("F10").Value = 815
$ws1.Range("F12").Value = 30
$ws1.Range("F13").Value = 570
$ws1.Range("F14").Value = 16
$ws1.Range("F15").Value = 18
$ws1.Range("F17").Value = 1768
$ws1.Range("F19").Value = 833
$ws1.Range("F21").Value = 189
$ws1.Range("F22").Value = 307
$ws1.Range("F23").Value = 513
$ws1.Range("F24").Value = 135
$ws1.Range("F28").Value = 2607
$ws1.Range("F29").Value = 173
$ws1.Range("F30").Value = 99
$ws1.Range("F31").Value = 55
$ws1.Range("F33").Value = 24
$ws1.Range("F34").Value = 267
$ws1.Range("F40").Value = 644
$ws1.Range("F41").Value = 81
$ws1.Range("F42").Value = 47
$ws1.Range("F44").Value = 60

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 233

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 233
$ws4.Range("F3").Value = 192
$ws4.Range("F7").Value = 5298
$ws4.Range("F8").Value = 446
$ws4.Range("F9").Value = 633
$ws4.Range("F12").Value = 914
$ws4.Range("F13").Value = 815
$ws4.Range("F17").Value = 30
$ws4.Range("F18").Value = 570
$ws4.Range("F19").Value = 16
$ws4.Range("F20").Value = 18
$ws4.Range("F23").Value = 1768
$ws4.Range("F25").Value = 833
$ws4.Range("F26").Value = 189
$ws4.Range("F27").Value = 307
$ws4.Range("F29").Value = 513
$ws4.Range("F30").Value = 135
$ws4.Range("F33").Value = 2607
$ws4.Range("F34").Value = 173
$ws4.Range("F35").Value = 99
$ws4.Range("F37").Value = 24
$ws4.Range("F38").Value = 267
$ws4.Range("F43").Value = 81
$ws4.Range("F44").Value = 47
$ws4.Range("F46").Value = 60
